$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "Existing ICDC Biobank filter" StatQuery: the Study-Files optional
# match was following the wrong edge (-->(s) study) - it should follow from
# the case (-->(c)), matching the per-tab queries' intent.
$fixedQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(c)
WITH DISTINCT f, sf,samp AS samp,c,demo,diag,s,p
WHERE s.clinical_study_designation IN ['COTC022'] and s.clinical_study_type in ['Clinical Trial'] and demo.neutered_indicator in ['No']
Optional MATCH(c)-->(r:registration)
WHERE r.registration_origin in ['Not Applicable']
RETURN
count(distinct p) AS Programs,
count(distinct s) AS Studies,
count(distinct c) AS Cases,
count(distinct samp) AS Samples,
count(distinct f) AS `Case Files`,
count(distinct sf) AS `Study Files`
'@

$ws.Range("C2").Value = $fixedQuery
$ws.Range("C3").Value = $fixedQuery
$ws.Range("C4").Value = $fixedQuery
$ws.Range("C5").Value = $fixedQuery

# Restore the user's on-screen selection/navigation state recorded at save time.
$ws.Range("J3").Select()
